$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: Site, Type, Latitude, Longitude, Elevation, Slope, Aspect, Compass
$ws.Range("A1").Value = "Site"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Latitude"
$ws.Range("D1").Value = "Longitude"
$ws.Range("H1").Value = "Compass"

# Data rows, reordered + new Site codes/Type/Compass columns added.
# Row 2: GOR (was CADCLIFFS row, now first)
$ws.Range("A2").Value = "GOR"
$ws.Range("B2").Value = "Cliff"
$ws.Range("C2").Value = 44.327796733333301
$ws.Range("D2").Value = -68.185402866666706
$ws.Range("E2").Value = 100.866666666667
$ws.Range("F2").Value = 31.866666666666699
$ws.Range("G2").Value = 151.46666666666701
$ws.Range("H2").Value = "SE"

# Row 3: SCT (was CAD row)
$ws.Range("A3").Value = "SCT"
$ws.Range("B3").Value = "Ledge"
$ws.Range("C3").Value = 44.3330482
$ws.Range("D3").Value = -68.223588000000007
$ws.Range("E3").Value = 930.73333333333301
$ws.Range("F3").Value = 9.4
$ws.Range("G3").Value = 138.13333333333301
$ws.Range("H3").Value = "SE"

# Row 4: STS (was STSAUV row)
$ws.Range("A4").Value = "STS"
$ws.Range("B4").Value = "Ledge"
$ws.Range("C4").Value = 44.310862800000002
$ws.Range("D4").Value = -68.326203933333304
$ws.Range("E4").Value = 597.13333333333298
$ws.Range("F4").Value = 13.3333333333333
$ws.Range("G4").Value = 262.60000000000002
$ws.Range("H4").Value = "W"

# Row 5: WON (was WOND row)
$ws.Range("A5").Value = "WON"
$ws.Range("B5").Value = "Flat"
$ws.Range("C5").Value = 44.237123400000002
$ws.Range("D5").Value = -68.316051799999997
$ws.Range("E5").Value = 53.4
$ws.Range("F5").Value = 3.7333333333333298
$ws.Range("G5").Value = 208.066666666667
$ws.Range("H5").Value = "S"

# Update selection to match the new edit focus
$ws.Range("A2:H5").Select()
